# Config.xlsx edit: "Set up cohort report and move email settings into config file"
#
# The Settings sheet gets 5 new rows inserted above the old "Framework"
# logging row (which slides from row 5 down to row 10). The new rows hold
# email-related configuration (previously hard-coded / elsewhere) that is
# now being moved into the config workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")
[void]$ws.Activate()

# Insert 5 blank rows above the current row 5 ("logF_BusinessProcessName").
# This pushes that row (and everything below it, including the long run of
# trailing blank formatted rows) down by 5, which is also what grows the
# sheet's used range from row 998 to row 1003.
[void]$ws.Rows("5:9").Insert()

# The newly-inserted rows don't inherit the sheet's standard explicit row
# height (14.25pt / custom) the way the rest of the sheet has it, so set it
# to match.
$ws.Rows("5:9").RowHeight = 14.25

# Row 5: Name of the Orchestrator asset holding the email credentials.
$ws.Range("A5").Value = "EmailCredentialAssetName"
$ws.Range("B5").Value = "Email Credentials"
$ws.Range("C5").Value = "Name of Orchestrator asset containing email credentials."

# Row 6: Orchestrator folder that contains the email credential asset.
$ws.Range("A6").Value = "EmailCredentialFolderPath"
$ws.Range("B6").Value = "P3 Automation"
$ws.Range("C6").Value = "Path to Orchestrator folder containing email credentials."

# Row 7: SMTP server used to send email.
$ws.Range("A7").Value = "EmailServer"
$ws.Range("B7").Value = "smtp.gmail.com"
$ws.Range("C7").Value = "Server for email origin."

# Row 8: SMTP port (numeric value).
$ws.Range("A8").Value = "EmailPort"
$ws.Range("B8").Value = 465
$ws.Range("B8").WrapText = $false
$ws.Range("C8").Value = "Port for email origin."

# Row 9: Override recipient address(es) for generated reports.
$ws.Range("A9").Value = "OverrideEmail"
$ws.Range("B9").Value = "wplee.327@gmail.com"
$ws.Range("C9").Value = "Override for email(s) that reports are sent to. Comma separated."

# The former row 5 content is now row 10 (untouched data, just shifted by
# the insert above) -- nothing else to change there.

# Move the visible selection to match the post-edit workbook (B11, just
# below the newly-entered settings block).
[void]$ws.Range("B11").Select()
